# Update countries & provincias Spain
#
# The upstream data refresh (scrape taken at 03:22 instead of 02:05) shuffled
# three country names within the shared-string table (Chequia, Botsuana and
# Belice each moved one slot earlier, ahead of Corea del Sur / Tanzania /
# Puerto Rico respectively) while every other row kept its position. Because
# column A in the worksheet is otherwise untouched, the net visible effect is
# that the country label shown on a handful of rows changes even though the
# row's statistics are independently refreshed too. We reproduce that net
# effect directly against the cells, plus the numeric "Casos totales / Nuevos
# casos / Casos activos / Recuperados / Casos criticos / Muertes hoy /
# Muertes" columns (B:H) that the new scrape updated, and the "Datos
# actualizados" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 03:22"

# --- Country labels that shifted because of the shared-string reorder ---
$ws.Range("A68").Value = "Chequia"
$ws.Range("A69").Value = "Corea del Sur"
$ws.Range("A70").Value = "Costa de Marfil"

$ws.Range("A158").Value = "Botsuana"
$ws.Range("A159").Value = "Tanzania"
$ws.Range("A160").Value = "Siria"
$ws.Range("A161").Value = "Taiwan"

$ws.Range("A195").Value = "Belice"
$ws.Range("A196").Value = "Puerto Rico"

# --- Refreshed statistics, row by row ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3694948
$ws.Range("C4").Value = 73311
$ws.Range("D4").Value = 1675823
$ws.Range("E4").Value = 1878008
$ws.Range("G4").Value = 962
$ws.Range("H4").Value = 141117

# Row 24: Canada
$ws.Range("E24").Value = 27601
$ws.Range("G24").Value = 17
$ws.Range("H24").Value = 8827

# Row 42: Panama
$ws.Range("B42").Value = 50373
$ws.Range("C42").Value = 1130
$ws.Range("D42").Value = 25842
$ws.Range("E42").Value = 23531
$ws.Range("G42").Value = 18
$ws.Range("H42").Value = 1000

# Row 68: Chequia
$ws.Range("C68").Value = 137
$ws.Range("D68").Value = 8640
$ws.Range("E68").Value = 4617
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 355

# Row 69: Corea del Sur
$ws.Range("B69").Value = 13612
$ws.Range("C69").Value = 61
$ws.Range("D69").Value = 12396
$ws.Range("E69").Value = 925
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 291

# Row 70: Costa de Marfil
$ws.Range("B70").Value = 13554
$ws.Range("C70").Value = 151
$ws.Range("D70").Value = 7363
$ws.Range("E70").Value = 6104
$ws.Range("H70").Value = 87

# Row 158: Botsuana
$ws.Range("B158").Value = 522
$ws.Range("C158").Value = 123
$ws.Range("D158").Value = 48
$ws.Range("E158").Value = 473
$ws.Range("H158").Value = 1

# Row 159: Tanzania
$ws.Range("B159").Value = 509
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 183
$ws.Range("E159").Value = 305
$ws.Range("H159").Value = 21

# Row 160: Siria
$ws.Range("B160").Value = 477
$ws.Range("C160").Value = 19
$ws.Range("D160").Value = 140
$ws.Range("E160").Value = 315
$ws.Range("H160").Value = 22

# Row 161: Taiwan
$ws.Range("B161").Value = 452
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 440
$ws.Range("E161").Value = 5
$ws.Range("H161").Value = 7

# Row 167: Guyana
$ws.Range("B167").Value = 315
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 140
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 19

# Row 178: Bermudas
$ws.Range("B178").Value = 152
$ws.Range("C178").Value = 2
$ws.Range("E178").Value = 6

# Row 184: Aruba
$ws.Range("B184").Value = 108
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 100
$ws.Range("E184").Value = 5

# Row 195: Belice
$ws.Range("B195").Value = 40
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 22
$ws.Range("E195").Value = 16

# Row 196: Puerto Rico
$ws.Range("D196").Value = 1
$ws.Range("E196").Value = 36
